$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 119.0815153333333
$ws.Range("H2").Value = 357.244546
$ws.Range("I2").Value = 0.431812569872284
$ws.Range("J2").Value = 0.4318125698722839
$ws.Range("M2").Value = 9.423852333333334
$ws.Range("N2").Value = 28.271557
$ws.Range("O2").Value = 0.06654336290212845
$ws.Range("P2").Value = 0.06654336290212845
$ws.Range("Q2").Value = 1122.206616130903
$ws.Range("R2").Value = 10099.85954517812
$ws.Range("S2").Value = 0.02873426054271209
$ws.Range("T2").Value = 0.02873426054271209
$ws.Range("G3").Value = 119.0815153333333
$ws.Range("H3").Value = 357.244546
$ws.Range("I3").Value = 0.431812569872284
$ws.Range("J3").Value = 0.4318125698722839
$ws.Range("O3").Value = 0.3572423751649123
$ws.Range("P3").Value = 0.3572423751649123
$ws.Range("Q3").Value = 6024.639265106325
$ws.Range("R3").Value = 54221.75338595692
$ws.Range("S3").Value = 0.1542617480872394
$ws.Range("T3").Value = 0.1542617480872394
$ws.Range("G4").Value = 119.0815153333333
$ws.Range("H4").Value = 357.244546
$ws.Range("I4").Value = 0.431812569872284
$ws.Range("J4").Value = 0.4318125698722839
$ws.Range("M4").Value = 26.84076266666667
$ws.Range("N4").Value = 80.522288
$ws.Range("O4").Value = 0.1895270158659356
$ws.Range("P4").Value = 0.1895270158659356
$ws.Range("Q4").Value = 3196.238691049028
$ws.Range("R4").Value = 28766.14821944125
$ws.Range("S4").Value = 0.08184014778129478
$ws.Range("T4").Value = 0.08184014778129478
$ws.Range("G5").Value = 119.0815153333333
$ws.Range("H5").Value = 357.244546
$ws.Range("I5").Value = 0.431812569872284
$ws.Range("J5").Value = 0.4318125698722839
$ws.Range("M5").Value = 54.762539
$ws.Range("N5").Value = 164.287617
$ws.Range("O5").Value = 0.3866872460670236
$ws.Range("P5").Value = 0.3866872460670236
$ws.Range("Q5").Value = 6521.206127620765
$ws.Range("R5").Value = 58690.85514858688
$ws.Range("S5").Value = 0.1669764134610377
$ws.Range("T5").Value = 0.1669764134610377
$ws.Range("I6").Value = 0.4460879372303943
$ws.Range("J6").Value = 0.4460879372303942
$ws.Range("M6").Value = 9.423852333333334
$ws.Range("N6").Value = 28.271557
$ws.Range("O6").Value = 0.06654336290212845
$ws.Range("P6").Value = 0.06654336290212845
$ws.Range("Q6").Value = 1159.305841152788
$ws.Range("R6").Value = 10433.75257037509
$ws.Range("S6").Value = 0.02968419149338403
$ws.Range("T6").Value = 0.02968419149338402
$ws.Range("I7").Value = 0.4460879372303943
$ws.Range("J7").Value = 0.4460879372303942
$ws.Range("O7").Value = 0.3572423751649123
$ws.Range("P7").Value = 0.3572423751649123
$ws.Range("S7").Value = 0.1593615142286024
$ws.Range("T7").Value = 0.1593615142286023
$ws.Range("I8").Value = 0.4460879372303943
$ws.Range("J8").Value = 0.4460879372303942
$ws.Range("M8").Value = 26.84076266666667
$ws.Range("N8").Value = 80.522288
$ws.Range("O8").Value = 0.1895270158659356
$ws.Range("P8").Value = 0.1895270158659356
$ws.Range("Q8").Value = 3301.903705600192
$ws.Range("R8").Value = 29717.13335040173
$ws.Range("S8").Value = 0.08454571555706741
$ws.Range("T8").Value = 0.08454571555706741
$ws.Range("I9").Value = 0.4460879372303943
$ws.Range("J9").Value = 0.4460879372303942
$ws.Range("M9").Value = 54.762539
$ws.Range("N9").Value = 164.287617
$ws.Range("O9").Value = 0.3866872460670236
$ws.Range("P9").Value = 0.3866872460670236
$ws.Range("Q9").Value = 6736.791822861828
$ws.Range("R9").Value = 60631.12640575645
$ws.Range("S9").Value = 0.1724965159513405
$ws.Range("T9").Value = 0.1724965159513405
$ws.Range("G10").Value = 33.50679633333333
$ws.Range("H10").Value = 100.520389
$ws.Range("I10").Value = 0.1215021138451521
$ws.Range("J10").Value = 0.121502113845152
$ws.Range("M10").Value = 9.423852333333334
$ws.Range("N10").Value = 28.271557
$ws.Range("O10").Value = 0.06654336290212845
$ws.Range("P10").Value = 0.06654336290212845
$ws.Range("Q10").Value = 315.7631008084081
$ws.Range("R10").Value = 2841.867907275673
$ws.Range("S10").Value = 0.008085159254973679
$ws.Range("T10").Value = 0.008085159254973677
$ws.Range("G11").Value = 33.50679633333333
$ws.Range("H11").Value = 100.520389
$ws.Range("I11").Value = 0.1215021138451521
$ws.Range("J11").Value = 0.121502113845152
$ws.Range("O11").Value = 0.3572423751649123
$ws.Range("P11").Value = 0.3572423751649123
$ws.Range("Q11").Value = 1695.194760266996
$ws.Range("R11").Value = 15256.75284240296
$ws.Range("S11").Value = 0.0434057037375997
$ws.Range("T11").Value = 0.04340570373759969
$ws.Range("G12").Value = 33.50679633333333
$ws.Range("H12").Value = 100.520389
$ws.Range("I12").Value = 0.1215021138451521
$ws.Range("J12").Value = 0.121502113845152
$ws.Range("M12").Value = 26.84076266666667
$ws.Range("N12").Value = 80.522288
$ws.Range("O12").Value = 0.1895270158659356
$ws.Range("P12").Value = 0.1895270158659356
$ws.Range("Q12").Value = 899.3479681033369
$ws.Range("R12").Value = 8094.131712930031
$ws.Range("S12").Value = 0.02302793305847485
$ws.Range("T12").Value = 0.02302793305847484
$ws.Range("G13").Value = 33.50679633333333
$ws.Range("H13").Value = 100.520389
$ws.Range("I13").Value = 0.1215021138451521
$ws.Range("J13").Value = 0.121502113845152
$ws.Range("M13").Value = 54.762539
$ws.Range("N13").Value = 164.287617
$ws.Range("O13").Value = 0.3866872460670236
$ws.Range("P13").Value = 0.3866872460670236
$ws.Range("Q13").Value = 1834.917240969223
$ws.Range("R13").Value = 16514.25516872301
$ws.Range("S13").Value = 0.04698331779410384
$ws.Range("T13").Value = 0.04698331779410383
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.16474
$ws.Range("H14").Value = 0.49422
$ws.Range("I14").Value = 0.000597379052169715
$ws.Range("J14").Value = 0.000597379052169715
$ws.Range("M14").Value = 9.423852333333334
$ws.Range("N14").Value = 28.271557
$ws.Range("O14").Value = 0.06654336290212845
$ws.Range("P14").Value = 0.06654336290212845
$ws.Range("Q14").Value = 1.552485433393334
$ws.Range("R14").Value = 13.97236890054
$ws.Range("S14").Value = 0.00003975161105865887
$ws.Range("T14").Value = 0.00003975161105865887
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.16474
$ws.Range("H15").Value = 0.49422
$ws.Range("I15").Value = 0.000597379052169715
$ws.Range("J15").Value = 0.000597379052169715
$ws.Range("O15").Value = 0.3572423751649123
$ws.Range("P15").Value = 0.3572423751649123
$ws.Range("Q15").Value = 8.334619103186666
$ws.Range("R15").Value = 75.01157192868
$ws.Range("S15").Value = 0.000213409111470873
$ws.Range("T15").Value = 0.000213409111470873
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.16474
$ws.Range("H16").Value = 0.49422
$ws.Range("I16").Value = 0.000597379052169715
$ws.Range("J16").Value = 0.000597379052169715
$ws.Range("M16").Value = 26.84076266666667
$ws.Range("N16").Value = 80.522288
$ws.Range("O16").Value = 0.1895270158659356
$ws.Range("P16").Value = 0.1895270158659356
$ws.Range("Q16").Value = 4.421747241706667
$ws.Range("R16").Value = 39.79572517536
$ws.Range("S16").Value = 0.0001132194690985471
$ws.Range("T16").Value = 0.0001132194690985471
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.16474
$ws.Range("H17").Value = 0.49422
$ws.Range("I17").Value = 0.000597379052169715
$ws.Range("J17").Value = 0.000597379052169715
$ws.Range("M17").Value = 54.762539
$ws.Range("N17").Value = 164.287617
$ws.Range("O17").Value = 0.3866872460670236
$ws.Range("P17").Value = 0.3866872460670236
$ws.Range("Q17").Value = 9.021580674859999
$ws.Range("R17").Value = 81.19422607373998
$ws.Range("S17").Value = 0.000230998860541636
$ws.Range("T17").Value = 0.000230998860541636
